# Hortaliza, Terminal La Palmera de La Serena - Albahaca
# Weekly price-report refresh: a new weekly observation is inserted as
# row 182 (Arica y Parinacota, $/paquete), pushing the previously-existing
# rows 182-203 down to 183-204 (their data is unchanged, just relocated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 182, shifting rows 182:203 down to
# 183:204 (and their formatting/styles, incl. the date number format) with them.
$ws.Rows.Item(182).Insert()

# Make the new row 182's date cell use the same date number format as its
# neighbours (style index 2 / "YYYY-MM-DD HH:MM:SS") instead of the default.
$ws.Range("D182").NumberFormat = $ws.Range("D183").NumberFormat

# Fill in the new weekly observation.
$ws.Range("A182").Value = 8
$ws.Range("B182").Value = "Terminal La Palmera de La Serena"
$ws.Range("C182").Value = "Coquimbo"
$ws.Range("D182").Value = 45131
$ws.Range("E182").Value = 4
$ws.Range("F182").Value = 100112052
$ws.Range("G182").Value = "Albahaca"
$ws.Range("H182").Value = "Sin especificar"
$ws.Range("I182").Value = "Primera"
$ws.Range("J182").Value = 1000
$ws.Range("K182").Value = 2500
$ws.Range("L182").Value = 3000
$ws.Range("M182").Value = 2750
$ws.Range("N182").Value = "$/paquete"
$ws.Range("O182").Value = "Región de Arica y Parinacota"
$ws.Range("P182").Value = 2750
$ws.Range("Q182").Value = 1
$ws.Range("R182").Value = "Hortaliza"
